# The presentation originally has 3 slides:
#   Slide 1 (sldId 256) - kept, unmodified
#   Slide 2 (sldId 258) - deleted
#   Slide 3 (sldId 257) - deleted
#
# Delete from the highest index down so earlier indices remain stable.
$p = $ppt.ActivePresentation

for ($i = $p.Slides.Count; $i -ge 2; $i--) {
    $p.Slides.Item($i).Delete()
}
